$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.285876989364624
$ws.Range("B1").Value = 3.843074798583984
$ws.Range("C1").Value = 5.939803123474121
$ws.Range("D1").Value = 2.539125919342041
$ws.Range("E1").Value = 1.86514937877655
